$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the alternating FAKE/REAL sender addresses in column A (rows 2-23)
#    so the pattern is shifted by one position (REAL now starts at row 2).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 23; $r++) {
  if ($r % 2 -eq 0) {
    $ws.Range("A$r").Value = "REAL@REAL.COM"
  } else {
    $ws.Range("A$r").Value = "FAKE@FAKE.COM"
  }
}

# ---------------------------------------------------------------------------
# 2. Rebuild the hyperlinks collection so that each A-cell's mailto link
#    matches its (new) displayed address, in the new relative order.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$order = @(3,4,5,7,9,11,13,15,17,19,21,23,6,8,10,12,14,16,18,20,22,2)
foreach ($r in $order) {
  if ($r % 2 -eq 0) {
    $addr = "mailto:REAL@REAL.COM"
  } else {
    $addr = "mailto:FAKE@FAKE.COM"
  }
  $ws.Hyperlinks.Add($ws.Range("A$r"), $addr)
}

# Re-apply the Hyperlink cell style (Hyperlinks.Add can otherwise create a
# duplicate style entry and bump the style index away from the original).
for ($r = 2; $r -le 23; $r++) {
  $ws.Range("A$r").Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 3. Add the new trailing blank row 24 (A24 keeps the Hyperlink-like style
#    but has no value/hyperlink of its own), extending the used range.
# ---------------------------------------------------------------------------
$ws.Range("A24").Style = "Hyperlink"
$ws.Range("A24").Value = ""

# ---------------------------------------------------------------------------
# 4. Update the view: selection moves to B25, top row scrolls to show row 18.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B25").Select()
